$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")

# Change the "id" column (A3:A6) from numeric ids to string ids matching
# the location's key used in the JSON (lowerCamelCase of the name).
$ws.Range("A3").Value = "bar"
$ws.Range("A4").Value = "livingRoom"
$ws.Range("A5").Value = "entrance"
$ws.Range("A6").Value = "bin"

# Update the view: active selection moved to A5 with A4 scrolled to the top.
$ws.Range("A5").Select()
